$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "22.408.09"
$ws.Range("E2").Value = "  +0.09%  "
Set-TextValue "D3" "1.571.82"
$ws.Range("E3").Value = "  +0.12%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue "D5" "1.002"
$ws.Range("E5").Value = "  +0.17%  "
Set-TextValue "D6" "291.55"
$ws.Range("E6").Value = "  +0.71%  "
Set-TextValue "D8" "49.88"
$ws.Range("E8").Value = "  +1.07%  "
Set-TextValue "D9" "0.3422"
$ws.Range("E9").Value = "  +1.32%  "
Set-TextValue "D10" "0.07648"
$ws.Range("E10").Value = "  +0.57%  "
Set-TextValue "D11" "1.151"
$ws.Range("E11").Value = "  -1.25%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  +0.16%  "
Set-TextValue "D13" "21.19"
$ws.Range("E13").Value = "  +0.01%  "
Set-TextValue "D14" "6.012"
$ws.Range("E14").Value = "  -0.71%  "
Set-TextValue "D15" "6.935"
$ws.Range("E15").Value = "  +0.58%  "
Set-TextValue "D16" "1.568.94"
$ws.Range("E16").Value = "  -0.15%  "
Set-TextValue "D17" "0.00001135"
$ws.Range("E17").Value = "  +0.10%  "
Set-TextValue "D18" "90.14"
$ws.Range("E18").Value = "  +1.27%  "
Set-TextValue "D19" "0.06731"
$ws.Range("E19").Value = "  -0.34%  "
Set-TextValue "D21" "16.76"
$ws.Range("E21").Value = "  +1.32%  "
Set-TextValue "D22" "6.226"
$ws.Range("E22").Value = "  -0.20%  "
Set-TextValue "D23" "12.03"
$ws.Range("E23").Value = "  +0.56%  "
Set-TextValue "D24" "22.405.04"
$ws.Range("E24").Value = "  -0.03%  "
Set-TextValue "D25" "2.398"
$ws.Range("E25").Value = "  +0.63%  "
Set-TextValue "D27" "20.21"
$ws.Range("E27").Value = "  +1.44%  "
Set-TextValue "D28" "147.06"
$ws.Range("E28").Value = "  +1.37%  "
Set-TextValue "D29" "5.026"
$ws.Range("E29").Value = "  +1.19%  "
Set-TextValue "D30" "126.86"
$ws.Range("E30").Value = "  +1.24%  "
Set-TextValue "D31" "1.749.99"
$ws.Range("E31").Value = "  +0.23%  "
Set-TextValue "D32" "6.171"
$ws.Range("E32").Value = "  -1.42%  "
Set-TextValue "D33" "2.007"
$ws.Range("E33").Value = "  +0.79%  "
Set-TextValue "D34" "0.9825"
$ws.Range("E34").Value = "  -5.44%  "
Set-TextValue "D36" "0.08478"
$ws.Range("E36").Value = "  +0.23%  "
Set-TextValue "D37" "0.02538"
$ws.Range("E37").Value = "  +0.19%  "
Set-TextValue "D38" "1.375"
$ws.Range("E38").Value = "  +10.46%  "
Set-TextValue "D39" "0.2315"
$ws.Range("E39").Value = "  -0.68%  "
Set-TextValue "D40" "0.06554"
$ws.Range("E40").Value = "  +0.67%  "
Set-TextValue "D41" "5.426"
$ws.Range("E41").Value = "  -2.24%  "
Set-TextValue "D42" "11.47"
$ws.Range("E42").Value = "  -2.95%  "
Set-TextValue "D43" "0.6393"
$ws.Range("E43").Value = "  +0.42%  "
Set-TextValue "D44" "14.14"
$ws.Range("E44").Value = "  -1.38%  "
Set-TextValue "D46" "3.804"
$ws.Range("E46").Value = "  +1.34%  "
Set-TextValue "D47" "0.5977"
$ws.Range("E47").Value = "  -0.13%  "
Set-TextValue "D50" "124.83"
$ws.Range("E50").Value = "  +1.12%  "
Set-TextValue "D51" "0.07332"
$ws.Range("E51").Value = "  +0.73%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E45").Value = "  +0.09%  "
Set-TextValue "D26" "2.660"

# Row 48/49: EOS and NEARProtocol swap places with new values
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D48" "2.097"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D49" "1.289"
$ws.Range("E49").Value = "  +2.34%  "
